# Apply the "Alexandria Convention" ontology-id formatting change:
# replace colon-separated ontology identifiers (PREFIX:ACCESSION) with
# underscore-separated ones (PREFIX_ACCESSION) in the example data row,
# and replace the formatted time value in Z2 with its literal text
# equivalent using the same underscore convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells F2, L2 use a "quote prefix" (force-text) cell style. Prefixing the
# new value with a leading apostrophe keeps that style/formatting intact
# (matches how Excel itself preserves the quote-prefix flag when a user
# types over such a cell) while still storing the clean text value.
$ws.Range("F2").Value = "'UBERON_0002107"
$ws.Range("G2").Value = "NCBITAXON_10090"
$ws.Range("I2").Value = "MGI_2159769"
$ws.Range("L2").Value = "'PATO_0000384"
$ws.Range("O2").Value = "UO_0000033"
$ws.Range("S2").Value = "UO_0000033"
$ws.Range("V2").Value = "UO_0000027"
$ws.Range("Y2").Value = "UO_0000187"
$ws.Range("AJ2").Value = "CHEBI_28119"

# Z2 held a numeric time-of-day value (formatted h:mm, displaying as a
# time like 12:12 PM). Replace it with the literal underscore-style text
# representation instead of a real time value.
$ws.Range("Z2").Value = "12_12_00 PM"

# Update the active selection/view: previously the sheet was scrolled so
# column BJ was left-most with BL2 selected; now the sheet is scrolled
# back to the left with L2 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("L2").Select()
